# Commit: "adding table format (row/column) to Workbook output"
#
# Every ObjTables data sheet carries a machine-generated header string in
# its first cell(s):
#   - A single workbook-level banner: "!!!ObjTables objTablesVersion='...' date='...'"
#     (only present once, on the very first sheet).
#   - A per-table banner on each sheet: "!!ObjTables type='Data' id='...' name='...' date='...' objTablesVersion='...'"
#
# This edit re-stamps the generation date/time on every one of those banner
# strings and appends a new `tableFormat='row'` attribute to each per-table
# banner (the workbook-level banner only gets the refreshed date).

$wb = $excel.ActiveWorkbook
$newDate = "2020-03-09 15:31:53"

function Update-ObjTablesCell($ws, $cellAddr) {
    $range = $ws.Range($cellAddr)
    $val = $range.Value2
    if ($val -eq $null) { return }
    $text = $val.ToString()
    if ($text.IndexOf("ObjTables") -ge 0 -and $text.IndexOf("date='") -ge 0) {
        $newText = [System.Text.RegularExpressions.Regex]::Replace($text, "date='[^']*'", "date='$newDate'")
        if ($newText.StartsWith("!!ObjTables") -and -not $newText.StartsWith("!!!ObjTables") -and $newText.IndexOf("tableFormat=") -lt 0) {
            $newText = $newText + " tableFormat='row'"
        }
        $range.Value2 = $newText
    }
}

foreach ($ws in $wb.Worksheets) {
    # Sheets ship protected (no password) so literal-string writes are
    # rejected unless we unprotect first; re-protect afterwards so the
    # sheet's lock state is unchanged from the user's point of view.
    $ws.Unprotect()
    Update-ObjTablesCell $ws "A1"
    Update-ObjTablesCell $ws "A2"
    $ws.Protect()
}
